$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 23:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 954360
$ws.Range("C4").Value = 29128
$ws.Range("E4").Value = 784296
$ws.Range("G4").Value = 1856
$ws.Range("H4").Value = 54049

# Row 16 - Canada
$ws.Range("B16").Value = 45138
$ws.Range("C16").Value = 1250
$ws.Range("D16").Value = 16351
$ws.Range("E16").Value = 26323
$ws.Range("G16").Value = 162
$ws.Range("H16").Value = 2464

# Row 89 - Guinea
$ws.Range("B89").Value = 996
$ws.Range("C89").Value = 42
$ws.Range("D89").Value = 208
$ws.Range("E89").Value = 781
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 7

# Row 94 - Principado de Andorra
$ws.Range("B94").Value = 738
$ws.Range("C94").Value = 7
$ws.Range("E94").Value = 354

# Row 132 - Ruanda
$ws.Range("B132").Value = 183
$ws.Range("C132").Value = 7
$ws.Range("D132").Value = 88
$ws.Range("E132").Value = 95
